$d = $word.ActiveDocument
$s = $d.Styles.Add("AbstractTitle", 1)
$s.BaseStyle = "Normal"
$s.NextParagraphStyle = "Abstract"
$s.QuickStyle = $true
$s.ParagraphFormat.KeepWithNext = $true
$s.ParagraphFormat.KeepTogether = $true
$s.ParagraphFormat.Alignment = 1
$s.ParagraphFormat.SpaceAfter = 0
$s.ParagraphFormat.SpaceBefore = 15
$s.Font.Size = 10
$s.Font.SizeBidirectional = 10
$s.Font.Bold = $true
$s.Font.Color = 9067060
